# "update report with US02/US03"
# Fill in actual size/time and mark the two Sprint1 backlog items
# (US02 row 6, US03 row 7) as completed, then leave the Sprint1 sheet
# active with I8 selected (matching the saved workbook view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# US02 (row 6): Act Size / Act Time / Completed
$ws.Range("G6").Value = 20
$ws.Range("H6").Value = 60
$ws.Range("I6").Value = "YES"

# US03 (row 7): Act Size / Act Time / Completed
$ws.Range("G7").Value = 20
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = "YES"

# Make Sprint1 the active sheet/tab and leave the selection where the
# author left it (I8) instead of the previous Team/B17 selections.
$ws.Activate()
$ws.Range("I8").Select()
